# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 262
$firstRow = 2

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45202
